# Update TPM-derived values in the LR-pairs sheet (C3-C3ar1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.04155
$ws.Range("H2").Value = 0.12465
$ws.Range("I2").Value = 0.0001466168179836329
$ws.Range("J2").Value = 0.0001466168179836329
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3939839999999999
$ws.Range("N2").Value = 1.181952
$ws.Range("O2").Value = 0.5362668594039028
$ws.Range("P2").Value = 0.5362668594039028
$ws.Range("Q2").Value = 0.0163700352
$ws.Range("R2").Value = 0.1473303168
$ws.Range("S2").Value = [double]"7.862574051587646E-05"
$ws.Range("T2").Value = [double]"7.862574051587646E-05"

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.04155
$ws.Range("H3").Value = 0.12465
$ws.Range("I3").Value = 0.0001466168179836329
$ws.Range("J3").Value = 0.0001466168179836329
$ws.Range("O3").Value = 0.4637331405960971
$ws.Range("P3").Value = 0.4637331405960971
$ws.Range("Q3").Value = 0.01415587725
$ws.Range("R3").Value = 0.12740289525
$ws.Range("S3").Value = [double]"6.799107746775639E-05"
$ws.Range("T3").Value = [double]"6.799107746775639E-05"

# Row 4
$ws.Range("I4").Value = 0.9992428949822291
$ws.Range("J4").Value = 0.9992428949822291
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3939839999999999
$ws.Range("N4").Value = 1.181952
$ws.Range("O4").Value = 0.5362668594039028
$ws.Range("P4").Value = 0.5362668594039028
$ws.Range("Q4").Value = 111.567292137216
$ws.Range("R4").Value = 1004.105629234944
$ws.Range("S4").Value = 0.5358608490737838
$ws.Range("T4").Value = 0.5358608490737838

# Row 5
$ws.Range("I5").Value = 0.9992428949822291
$ws.Range("J5").Value = 0.9992428949822291
$ws.Range("O5").Value = 0.4637331405960971
$ws.Range("P5").Value = 0.4637331405960971
$ws.Range("S5").Value = 0.4633820459084451
$ws.Range("T5").Value = 0.4633820459084451

# Row 6
$ws.Range("I6").Value = [double]"0.0006104881997874136"
$ws.Range("J6").Value = [double]"0.0006104881997874135"
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3939839999999999
$ws.Range("N6").Value = 1.181952
$ws.Range("O6").Value = 0.5362668594039028
$ws.Range("P6").Value = 0.5362668594039028
$ws.Range("Q6").Value = 0.06816212121599999
$ws.Range("R6").Value = 0.613459090944
$ws.Range("S6").Value = [double]"0.0003273845896031387"
$ws.Range("T6").Value = [double]"0.0003273845896031386"

# Row 7
$ws.Range("I7").Value = [double]"0.0006104881997874136"
$ws.Range("J7").Value = [double]"0.0006104881997874135"
$ws.Range("O7").Value = 0.4637331405960971
$ws.Range("P7").Value = 0.4637331405960971
$ws.Range("S7").Value = [double]"0.0002831036101842749"
$ws.Range("T7").Value = [double]"0.0002831036101842748"
